$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P1 and Q1, continuing the sequence (values 14 and 15),
# using the same style as the rest of the header row (style index 1 -> bold,
# bordered, centered). Copy the format from O1 so no new cell styles get
# created in styles.xml (keeps cellXfs count identical to the source file).
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = 0

# For data rows 2-25: swap values in columns I, K, M, O and add new columns P, Q
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: was 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: was 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: was 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: was 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new column, value 2
    $ws.Cells.Item($r, 17).Value = 2  # Q: new column, value 2
}
